$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column D: FHIR element-path annotations.
# Values are entered in the order that first introduces each new distinct
# string, so the shared-string table gets built up the same way.
$ws.Range("D2").Value = ".code"
$ws.Range("D4").Value = ".category"
$ws.Range("D13").Value = ".alias"
$ws.Range("D24").Value = ".mapping"
$ws.Range("D26").Value = ".name"
$ws.Range("D12").Value = ".status"
$ws.Range("D19").Value = ".question"
$ws.Range("D11").Value = ".comments"
$ws.Range("D18").Value = ".binding.reference"
$ws.Range("D21").Value = "(how is this different from alias?)"

# Remaining rows that repeat the ".category" annotation
$ws.Range("D5").Value = ".category"
$ws.Range("D6").Value = ".category"
$ws.Range("D7").Value = ".category"
$ws.Range("D8").Value = ".category"
$ws.Range("D9").Value = ".category"

# Move the active cell/selection to D3, matching the saved view state
$ws.Range("D3").Select()
